$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sFcqD252"
$ws.Range("B2").Value = 231006272
$ws.Range("C2").Value = "mysqexo49"
$ws.Range("D2").Value = "psSX#8&2"
$ws.Range("F2").Value = "WEJzhDQR"
$ws.Range("G2").Value = "jsmy"
